# Auto-generated edit script: update Phantom_Profits market-price figures
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(9, 8).Value = 624.5  # H9: 838.2857 -> 624.5
$ws.Cells.Item(9, 9).Value = 672.7778  # I9: 946.3333 -> 672.7778
$ws.Cells.Item(9, 11).Value = 672.7778  # K9: 946.3333 -> 672.7778
$ws.Cells.Item(9, 13).Value = -503.7778  # M9: -777.3333 -> -503.7778

$ws.Cells.Item(32, 8).Value = 6369.857  # H32: 6831.5 -> 6369.857
$ws.Cells.Item(32, 9).Value = 5533  # I32: 5833 -> 5533
$ws.Cells.Item(32, 10).Value = 6997.5  # J32: 7830 -> 6997.5
$ws.Cells.Item(32, 11).Value = 5533  # K32: 5833 -> 5533
$ws.Cells.Item(32, 12).Value = 6997.5  # L32: 7830 -> 6997.5
$ws.Cells.Item(32, 13).Value = -5207  # M32: -5507 -> -5207
$ws.Cells.Item(32, 14).Value = -7649.5  # N32: -8482 -> -7649.5

$ws.Cells.Item(100, 8).Value = 2267.375  # H100: 2412.2727 -> 2267.375
$ws.Cells.Item(100, 9).Value = 2333.4348  # I100: 2412.2727 -> 2333.4348
$ws.Cells.Item(100, 10).Value = 748  # J100: 0 -> 748
$ws.Cells.Item(100, 11).Value = 2333.4348  # K100: 2412.2727 -> 2333.4348
$ws.Cells.Item(100, 12).Value = 748  # L100: 0 -> 748
$ws.Cells.Item(100, 13).Value = -1792.4348  # M100: -1871.2727 -> -1792.4348
$ws.Cells.Item(100, 14).Value = -1830  # N100: None -> -1830

$ws.Cells.Item(125, 8).Value = 981.65216  # H125: 799 -> 981.65216
$ws.Cells.Item(125, 10).Value = 981.65216  # J125: 799 -> 981.65216
$ws.Cells.Item(125, 12).Value = 8834.86944  # L125: 7191 -> 8834.86944
$ws.Cells.Item(125, 14).Value = -13754.86944  # N125: -12111 -> -13754.86944

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2529.4722  # H32: 2576.7222 -> 2529.4722
$ws.Cells.Item(32, 9).Value = 2379.0857  # I32: 2427.6858 -> 2379.0857
$ws.Cells.Item(32, 11).Value = 2379.0857  # K32: 2427.6858 -> 2379.0857
$ws.Cells.Item(32, 13).Value = -2092.0857  # M32: -2140.6858 -> -2092.0857

$ws.Cells.Item(45, 8).Value = 2061.7856  # H45: 2253.5 -> 2061.7856
$ws.Cells.Item(45, 9).Value = 1651.1538  # I45: 1785.6364 -> 1651.1538
$ws.Cells.Item(45, 11).Value = 1651.1538  # K45: 1785.6364 -> 1651.1538
$ws.Cells.Item(45, 13).Value = -1274.1538  # M45: -1408.6364 -> -1274.1538

$ws.Cells.Item(61, 8).Value = 2148.3  # H61: 2222.875 -> 2148.3
$ws.Cells.Item(61, 9).Value = 1783.4286  # I61: 1797.3334 -> 1783.4286
$ws.Cells.Item(61, 10).Value = 2999.6667  # J61: 3499.5 -> 2999.6667
$ws.Cells.Item(61, 11).Value = 1783.4286  # K61: 1797.3334 -> 1783.4286
$ws.Cells.Item(61, 12).Value = 2999.6667  # L61: 3499.5 -> 2999.6667
$ws.Cells.Item(61, 13).Value = -1571.4286  # M61: -1585.3334 -> -1571.4286
$ws.Cells.Item(61, 14).Value = -3423.6667  # N61: -3923.5 -> -3423.6667

$ws.Cells.Item(74, 8).Value = 3760.7144  # H74: 3973.4285 -> 3760.7144
$ws.Cells.Item(74, 9).Value = 3425.4546  # I74: 3566.9 -> 3425.4546
$ws.Cells.Item(74, 10).Value = 4990  # J74: 4989.75 -> 4990
$ws.Cells.Item(74, 11).Value = 3425.4546  # K74: 3566.9 -> 3425.4546
$ws.Cells.Item(74, 12).Value = 4990  # L74: 4989.75 -> 4990
$ws.Cells.Item(74, 13).Value = -2551.4546  # M74: -2692.9 -> -2551.4546
$ws.Cells.Item(74, 14).Value = -6738  # N74: -6737.75 -> -6738

$ws.Cells.Item(77, 8).Value = 3760.7144  # H77: 3973.4285 -> 3760.7144
$ws.Cells.Item(77, 9).Value = 3425.4546  # I77: 3566.9 -> 3425.4546
$ws.Cells.Item(77, 10).Value = 4990  # J77: 4989.75 -> 4990
$ws.Cells.Item(77, 11).Value = 17127.273  # K77: 17834.5 -> 17127.273
$ws.Cells.Item(77, 12).Value = 24950  # L77: 24948.75 -> 24950
$ws.Cells.Item(77, 13).Value = -12759.273  # M77: -13466.5 -> -12759.273
$ws.Cells.Item(77, 14).Value = -33686  # N77: -33684.75 -> -33686

$ws.Cells.Item(102, 8).Value = 4235.4  # H102: 4235.6 -> 4235.4
$ws.Cells.Item(102, 9).Value = 3419.25  # I102: 3419.5 -> 3419.25
$ws.Cells.Item(102, 11).Value = 3419.25  # K102: 3419.5 -> 3419.25
$ws.Cells.Item(102, 13).Value = -1797.25  # M102: -1797.5 -> -1797.25

$ws.Cells.Item(110, 8).Value = 22943.75  # H110: 24505.428 -> 22943.75
$ws.Cells.Item(110, 9).Value = 20381  # I110: 23170.666 -> 20381
$ws.Cells.Item(110, 11).Value = 20381  # K110: 23170.666 -> 20381
$ws.Cells.Item(110, 13).Value = -18336  # M110: -21125.666 -> -18336

$ws.Cells.Item(136, 8).Value = 2148.3  # H136: 2222.875 -> 2148.3
$ws.Cells.Item(136, 9).Value = 1783.4286  # I136: 1797.3334 -> 1783.4286
$ws.Cells.Item(136, 10).Value = 2999.6667  # J136: 3499.5 -> 2999.6667
$ws.Cells.Item(136, 11).Value = 5350.2858  # K136: 5392.0002 -> 5350.2858
$ws.Cells.Item(136, 12).Value = 8999.000100000001  # L136: 10498.5 -> 8999.000100000001
$ws.Cells.Item(136, 13).Value = -2800.2858  # M136: -2842.0002 -> -2800.2858
$ws.Cells.Item(136, 14).Value = -14099.0001  # N136: -15598.5 -> -14099.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 4466.722  # H80: 4521.1055 -> 4466.722
$ws.Cells.Item(80, 10).Value = 18474.5  # J80: 15879.6 -> 18474.5
$ws.Cells.Item(80, 12).Value = 18474.5  # L80: 15879.6 -> 18474.5
$ws.Cells.Item(80, 14).Value = -20470.5  # N80: -17875.6 -> -20470.5

$ws.Cells.Item(83, 8).Value = 4466.722  # H83: 4521.1055 -> 4466.722
$ws.Cells.Item(83, 10).Value = 18474.5  # J83: 15879.6 -> 18474.5
$ws.Cells.Item(83, 12).Value = 92372.5  # L83: 79398 -> 92372.5
$ws.Cells.Item(83, 14).Value = -102356.5  # N83: -89382 -> -102356.5

$ws.Cells.Item(134, 8).Value = 1751.4  # H134: 1776 -> 1751.4
$ws.Cells.Item(134, 9).Value = 1751.4  # I134: 1776 -> 1751.4
$ws.Cells.Item(134, 11).Value = 5254.200000000001  # K134: 5328 -> 5254.200000000001
$ws.Cells.Item(134, 13).Value = -2719.200000000001  # M134: -2793 -> -2719.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(20, 8).Value = 84997  # H20: 90000 -> 84997
$ws.Cells.Item(20, 10).Value = 84997  # J20: 90000 -> 84997
$ws.Cells.Item(20, 12).Value = 84997  # L20: 90000 -> 84997
$ws.Cells.Item(20, 14).Value = -85469  # N20: -90472 -> -85469

$ws.Cells.Item(22, 8).Value = 2331.75  # H22: 7500 -> 2331.75
$ws.Cells.Item(22, 9).Value = 413.5  # I22: 0 -> 413.5
$ws.Cells.Item(22, 10).Value = 4250  # J22: 7500 -> 4250
$ws.Cells.Item(22, 11).Value = 413.5  # K22: 0 -> 413.5
$ws.Cells.Item(22, 12).Value = 4250  # L22: 7500 -> 4250
$ws.Cells.Item(22, 13).Value = -63.5  # M22: None -> -63.5
$ws.Cells.Item(22, 14).Value = -4950  # N22: -8200 -> -4950

$ws.Cells.Item(30, 8).Value = 84997  # H30: 90000 -> 84997
$ws.Cells.Item(30, 10).Value = 84997  # J30: 90000 -> 84997
$ws.Cells.Item(30, 12).Value = 84997  # L30: 90000 -> 84997
$ws.Cells.Item(30, 14).Value = -85179  # N30: -90182 -> -85179

$ws.Cells.Item(31, 8).Value = 2346.4285  # H31: 3681.7144 -> 2346.4285
$ws.Cells.Item(31, 9).Value = 3422.2  # I31: 4027.75 -> 3422.2
$ws.Cells.Item(31, 10).Value = 1748.7778  # J31: 3220.3333 -> 1748.7778
$ws.Cells.Item(31, 11).Value = 3422.2  # K31: 4027.75 -> 3422.2
$ws.Cells.Item(31, 12).Value = 1748.7778  # L31: 3220.3333 -> 1748.7778
$ws.Cells.Item(31, 13).Value = -3127.2  # M31: -3732.75 -> -3127.2
$ws.Cells.Item(31, 14).Value = -2338.7778  # N31: -3810.3333 -> -2338.7778

$ws.Cells.Item(34, 8).Value = 2346.4285  # H34: 3681.7144 -> 2346.4285
$ws.Cells.Item(34, 9).Value = 3422.2  # I34: 4027.75 -> 3422.2
$ws.Cells.Item(34, 10).Value = 1748.7778  # J34: 3220.3333 -> 1748.7778
$ws.Cells.Item(34, 11).Value = 3422.2  # K34: 4027.75 -> 3422.2
$ws.Cells.Item(34, 12).Value = 1748.7778  # L34: 3220.3333 -> 1748.7778
$ws.Cells.Item(34, 13).Value = -3220.2  # M34: -3825.75 -> -3220.2
$ws.Cells.Item(34, 14).Value = -2152.7778  # N34: -3624.3333 -> -2152.7778

$ws.Cells.Item(47, 8).Value = 25399.6  # H47: 24349.666 -> 25399.6
$ws.Cells.Item(47, 9).Value = 21999  # I47: 20549.5 -> 21999
$ws.Cells.Item(47, 11).Value = 21999  # K47: 20549.5 -> 21999
$ws.Cells.Item(47, 13).Value = -21433  # M47: -19983.5 -> -21433

$ws.Cells.Item(94, 8).Value = 1006  # H94: 1061.5 -> 1006
$ws.Cells.Item(94, 9).Value = 1006  # I94: 1061.5 -> 1006
$ws.Cells.Item(94, 11).Value = 1006  # K94: 1061.5 -> 1006
$ws.Cells.Item(94, 13).Value = -555  # M94: -610.5 -> -555

$ws.Cells.Item(99, 8).Value = 1875.3334  # H99: 2206 -> 1875.3334
$ws.Cells.Item(99, 10).Value = 1214  # J99: 0 -> 1214
$ws.Cells.Item(99, 12).Value = 1214  # L99: 0 -> 1214
$ws.Cells.Item(99, 14).Value = -4210  # N99: None -> -4210

$ws.Cells.Item(122, 8).Value = 1954.2307  # H122: 2089.1667 -> 1954.2307
$ws.Cells.Item(122, 9).Value = 1976.125  # I122: 2210.5715 -> 1976.125
$ws.Cells.Item(122, 11).Value = 5928.375  # K122: 6631.7145 -> 5928.375
$ws.Cells.Item(122, 13).Value = -3478.375  # M122: -4181.7145 -> -3478.375

$ws.Cells.Item(126, 8).Value = 1875.3334  # H126: 2206 -> 1875.3334
$ws.Cells.Item(126, 10).Value = 1214  # J126: 0 -> 1214
$ws.Cells.Item(126, 12).Value = 3642  # L126: 0 -> 3642
$ws.Cells.Item(126, 14).Value = -8582  # N126: None -> -8582

$ws.Cells.Item(128, 8).Value = 84997  # H128: 90000 -> 84997
$ws.Cells.Item(128, 10).Value = 84997  # J128: 90000 -> 84997
$ws.Cells.Item(128, 12).Value = 84997  # L128: 90000 -> 84997
$ws.Cells.Item(128, 14).Value = -94957  # N128: -99960 -> -94957

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 1052.283  # H113: 1053.8654 -> 1052.283
$ws.Cells.Item(113, 9).Value = 1059.6875  # I113: 1059.75 -> 1059.6875
$ws.Cells.Item(113, 10).Value = 981.2  # J113: 983.25 -> 981.2
$ws.Cells.Item(113, 11).Value = 3179.0625  # K113: 3179.25 -> 3179.0625
$ws.Cells.Item(113, 12).Value = 2943.6  # L113: 2949.75 -> 2943.6
$ws.Cells.Item(113, 13).Value = -1009.0625  # M113: -1009.25 -> -1009.0625
$ws.Cells.Item(113, 14).Value = -7283.6  # N113: -7289.75 -> -7283.6

$ws.Cells.Item(131, 8).Value = 1886  # H131: 1743.3334 -> 1886
$ws.Cells.Item(131, 9).Value = 1176.6666  # I131: 1140 -> 1176.6666
$ws.Cells.Item(131, 11).Value = 3529.9998  # K131: 3420 -> 3529.9998
$ws.Cells.Item(131, 13).Value = 1510.0002  # M131: 1620 -> 1510.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 8262.857  # H43: 6501 -> 8262.857
$ws.Cells.Item(43, 9).Value = 560.17645  # I43: 640.1667 -> 560.17645
$ws.Cells.Item(43, 10).Value = 40999.25  # J43: 41666 -> 40999.25
$ws.Cells.Item(43, 11).Value = 560.17645  # K43: 640.1667 -> 560.17645
$ws.Cells.Item(43, 12).Value = 40999.25  # L43: 41666 -> 40999.25
$ws.Cells.Item(43, 13).Value = -409.17645  # M43: -489.1667 -> -409.17645
$ws.Cells.Item(43, 14).Value = -41301.25  # N43: -41968 -> -41301.25

$ws.Cells.Item(122, 8).Value = 3217.5  # H122: 3180.375 -> 3217.5
$ws.Cells.Item(122, 9).Value = 3217.5  # I122: 3308.8333 -> 3217.5
$ws.Cells.Item(122, 10).Value = 0  # J122: 2795 -> 0
$ws.Cells.Item(122, 11).Value = 9652.5  # K122: 9926.499899999999 -> 9652.5
$ws.Cells.Item(122, 12).Value = 0  # L122: 8385 -> 0
$ws.Cells.Item(122, 13).Value = -7202.5  # M122: -7476.499899999999 -> -7202.5
$ws.Cells.Item(122, 14).ClearContents()  # N122: -13285 -> (removed)

$ws.Cells.Item(126, 8).Value = 1999  # H126: 1999.5 -> 1999
$ws.Cells.Item(126, 9).Value = 1998.5  # I126: 1999 -> 1998.5
$ws.Cells.Item(126, 11).Value = 5995.5  # K126: 5997 -> 5995.5
$ws.Cells.Item(126, 13).Value = -3525.5  # M126: -3527 -> -3525.5

$ws.Cells.Item(128, 8).Value = 94900  # H128: 94950 -> 94900
$ws.Cells.Item(128, 9).Value = 94900  # I128: 0 -> 94900
$ws.Cells.Item(128, 10).Value = 0  # J128: 94950 -> 0
$ws.Cells.Item(128, 11).Value = 94900  # K128: 0 -> 94900
$ws.Cells.Item(128, 12).Value = 0  # L128: 94950 -> 0
$ws.Cells.Item(128, 13).Value = -89920  # M128: None -> -89920
$ws.Cells.Item(128, 14).ClearContents()  # N128: -104910 -> (removed)

$ws.Cells.Item(132, 8).Value = 4425  # H132: 4364.143 -> 4425
$ws.Cells.Item(132, 10).Value = 0  # J132: 3999 -> 0
$ws.Cells.Item(132, 12).Value = 0  # L132: 11997 -> 0
$ws.Cells.Item(132, 14).ClearContents()  # N132: -17057 -> (removed)

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 245  # H93: 244.66667 -> 245
$ws.Cells.Item(93, 9).Value = 245  # I93: 244.66667 -> 245
$ws.Cells.Item(93, 11).Value = 245  # K93: 244.66667 -> 245
$ws.Cells.Item(93, 13).Value = 1003  # M93: 1003.33333 -> 1003

$ws.Cells.Item(122, 8).Value = 7749.5  # H122: 8199.200000000001 -> 7749.5
$ws.Cells.Item(122, 10).Value = 8999.5  # J122: 9332.333000000001 -> 8999.5
$ws.Cells.Item(122, 12).Value = 26998.5  # L122: 27996.999 -> 26998.5
$ws.Cells.Item(122, 14).Value = -31898.5  # N122: -32896.999 -> -31898.5

$ws.Cells.Item(136, 8).Value = 2201.3635  # H136: 2268.889 -> 2201.3635
$ws.Cells.Item(136, 9).Value = 2182.5  # I136: 2253.75 -> 2182.5
$ws.Cells.Item(136, 11).Value = 6547.5  # K136: 6761.25 -> 6547.5
$ws.Cells.Item(136, 13).Value = -3997.5  # M136: -4211.25 -> -3997.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(28, 8).Value = 0  # H28: 8500 -> 0
$ws.Cells.Item(28, 9).Value = 0  # I28: 8500 -> 0
$ws.Cells.Item(28, 11).Value = 0  # K28: 8500 -> 0
$ws.Cells.Item(28, 13).ClearContents()  # M28: -8152 -> (removed)

$ws.Cells.Item(54, 8).Value = 15687.857  # H54: 16270.228 -> 15687.857
$ws.Cells.Item(54, 10).Value = 29611.75  # J54: 29389.4 -> 29611.75
$ws.Cells.Item(54, 12).Value = 29611.75  # L54: 29389.4 -> 29611.75
$ws.Cells.Item(54, 14).Value = -30651.75  # N54: -30429.4 -> -30651.75

$ws.Cells.Item(81, 8).Value = 2826.1538  # H81: 3399.8 -> 2826.1538
$ws.Cells.Item(81, 9).Value = 2957.6667  # I81: 4249.5 -> 2957.6667
$ws.Cells.Item(81, 10).Value = 2713.4285  # J81: 2833.3333 -> 2713.4285
$ws.Cells.Item(81, 11).Value = 5915.3334  # K81: 8499 -> 5915.3334
$ws.Cells.Item(81, 12).Value = 5426.857  # L81: 5666.6666 -> 5426.857
$ws.Cells.Item(81, 13).Value = -4854.3334  # M81: -7438 -> -4854.3334
$ws.Cells.Item(81, 14).Value = -7548.857  # N81: -7788.6666 -> -7548.857

$ws.Cells.Item(84, 8).Value = 2826.1538  # H84: 3399.8 -> 2826.1538
$ws.Cells.Item(84, 9).Value = 2957.6667  # I84: 4249.5 -> 2957.6667
$ws.Cells.Item(84, 10).Value = 2713.4285  # J84: 2833.3333 -> 2713.4285
$ws.Cells.Item(84, 11).Value = 29576.667  # K84: 42495 -> 29576.667
$ws.Cells.Item(84, 12).Value = 27134.285  # L84: 28333.333 -> 27134.285
$ws.Cells.Item(84, 13).Value = -24272.667  # M84: -37191 -> -24272.667
$ws.Cells.Item(84, 14).Value = -37742.285  # N84: -38941.333 -> -37742.285

$ws.Cells.Item(113, 8).Value = 563.6667  # H113: 571.8889 -> 563.6667
$ws.Cells.Item(113, 9).Value = 447.5  # I113: 466 -> 447.5
$ws.Cells.Item(113, 11).Value = 1342.5  # K113: 1398 -> 1342.5
$ws.Cells.Item(113, 13).Value = 827.5  # M113: 772 -> 827.5

$ws.Cells.Item(122, 8).Value = 775.375  # H122: 1438 -> 775.375
$ws.Cells.Item(122, 9).Value = 775.375  # I122: 1438 -> 775.375
$ws.Cells.Item(122, 11).Value = 2326.125  # K122: 4314 -> 2326.125
$ws.Cells.Item(122, 13).Value = 123.875  # M122: -1864 -> 123.875

$ws.Cells.Item(136, 8).Value = 11684.074  # H136: 11412.607 -> 11684.074
$ws.Cells.Item(136, 9).Value = 10482.542  # I136: 10226.56 -> 10482.542
$ws.Cells.Item(136, 11).Value = 31447.626  # K136: 30679.68 -> 31447.626
$ws.Cells.Item(136, 13).Value = -28897.626  # M136: -28129.68 -> -28897.626

Write-Output "Applied 217 cell updates across 8 sheets."
